$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.731.37'
$ws.Range('E2').Value = '  -6.77%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.695.04'
$ws.Range('E3').Value = '  -5.98%  '
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.29'
$ws.Range('E5').Value = '  -5.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5012'
$ws.Range('E6').Value = '  -15.84%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2590'
$ws.Range('E8').Value = '  -6.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '21.74'
$ws.Range('E9').Value = '  -6.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06104'
$ws.Range('E10').Value = '  -10.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07305'
$ws.Range('E11').Value = '  -2.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.689.36'
$ws.Range('E12').Value = '  -6.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.422'
$ws.Range('E13').Value = '  -6.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.926.27'
$ws.Range('E14').Value = '  -5.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5712'
$ws.Range('E15').Value = '  -8.68%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008174'
$ws.Range('E16').Value = '  -11.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.89'
$ws.Range('E17').Value = '  -13.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.734.66'
$ws.Range('E18').Value = '  -6.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.010'
$ws.Range('E19').Value = '  -8.11%  '
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.73'
$ws.Range('E21').Value = '  -6.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '183.72'
$ws.Range('E22').Value = '  -12.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.201'
$ws.Range('E23').Value = '  -9.25%  '
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.26'
$ws.Range('E25').Value = '  -5.95%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.603'
$ws.Range('E26').Value = '  -2.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1138'
$ws.Range('E28').Value = '  -6.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.315'
$ws.Range('E29').Value = '  -9.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05592'
$ws.Range('E30').Value = '  -10.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.327'
$ws.Range('E31').Value = '  -6.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.454'
$ws.Range('E32').Value = '  -8.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.441'
$ws.Range('E33').Value = '  -7.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.653'
$ws.Range('E34').Value = '  -3.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.006'
$ws.Range('E35').Value = '  -4.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.407'
$ws.Range('E36').Value = '  -3.78%  '
$ws.Range('E37').Value = '  -7.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.631'
$ws.Range('E38').Value = '  -3.30%  '
$ws.Range('E39').Value = '  -7.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.069.40'
$ws.Range('E40').Value = '  -5.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.876'
$ws.Range('E41').Value = '  -8.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8510'
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '98.56'
$ws.Range('E44').Value = '  -2.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.856.30'
$ws.Range('E45').Value = '  -5.35%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000107'
$ws.Range('E46').Value = '  -3.91%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.25'
$ws.Range('E47').Value = '  -6.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.006'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.082'
$ws.Range('E49').Value = '  -2.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4334'
$ws.Range('E50').Value = '  -3.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05207'
$ws.Range('E51').Value = '  -4.57%  '
